$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fixed the Cases-tab Cypher query: removed the erroneous trailing `co`/Cohort
# variable & RETURN column (the OPTIONAL MATCH (co:cohort) clause stays, but
# the query no longer returns a Cohort column that callers didn't expect).
$caseQuery = 'MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed IN [''Mixed Breed'']  MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '''') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '''') AS `Study Code` ,
        coalesce(s.clinical_study_type, '''') AS  `Study Type`,
        coalesce(demo.breed, '''') AS Breed ,
        coalesce(diag.disease_term, '''') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '''') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '''') AS Age ,
        coalesce(demo.sex, '''') AS Sex ,
        coalesce(demo.neutered_indicator, '''') AS `Neutered Status`,
        coalesce(demo.weight, '''') AS `Weight (kg)`,
        coalesce(diag.best_response, '''') AS `Response to Treatment`'
$ws.Range("B2").Value = $caseQuery

# Match author's final selection/view state.
[void]$ws.Range("B2").Select()

# Re-fit the row heights for the three wrapped query cells now that the text
# has changed length (Excel recalculates these automatically when content
# changes; we pin the values to match).
$ws.Rows.Item(2).RowHeight = 230.4
$ws.Rows.Item(3).RowHeight = 230.4
$ws.Rows.Item(4).RowHeight = 216
